# Applies: "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Replaces the worker detail rows (16-78) so that:
#  - Rows 16-76 hold DINA PATRICIA HERNANDEZ TEJEDOR (CC 45760836) with her 61 mora periods,
#    descending from 2309 down to 1809 (row 16 keeps the special 18750 Valor Mora for period 2309).
#  - Rows 77-78 hold YEIMI ARELLANO RODRIGUEZ (CC 1050950910) with periods 1810 and 1809.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$dinaDoc    = "45760836"
$dinaName   = "DINA PATRICIA HERNANDEZ TEJEDOR"
$yeimiDoc   = "1050950910"
$yeimiName  = "YEIMI ARELLANO RODRIGUEZ"

$dinaPeriods = @(
    "2309","2308","2307","2306","2305","2304","2303","2302","2301",
    "2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201",
    "2112","2111","2110","2109","2108","2107","2106","2105","2104","2103","2102","2101",
    "2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809"
)

$row = 16
foreach ($periodo in $dinaPeriods) {
    $ws.Cells.Item($row, 3).Value = $dinaDoc
    $ws.Cells.Item($row, 4).Value = $dinaName
    $ws.Cells.Item($row, 5).Value = $periodo
    if ($row -eq 16) {
        $ws.Cells.Item($row, 6).Value = 18750
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }
    $row = $row + 1
}

$yeimiPeriods = @("1810", "1809")
foreach ($periodo in $yeimiPeriods) {
    $ws.Cells.Item($row, 3).Value = $yeimiDoc
    $ws.Cells.Item($row, 4).Value = $yeimiName
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = 31249
    $row = $row + 1
}
